$d = $word.ActiveDocument

# Change 1: "refresh the page by clicking the " -> "refresh the page by clicking "
$d.Content.Find.Execute(
    "refresh the page by clicking the ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "refresh the page by clicking ",
    2
)

Write-Host "done"
